$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-order the "Bugs" sheet so it sits right after "Tasks" (position 2).
# ---------------------------------------------------------------------------
$bugsSheet = $wb.Worksheets.Item("Bugs")
$bugsSheet.Move($wb.Worksheets.Item("Generate Salary screen"))

# ---------------------------------------------------------------------------
# 2. Add the new task entries to the "Tasks" sheet.
# ---------------------------------------------------------------------------
$tasks = $wb.Worksheets.Item("Tasks")

# Row 47 gains a start / end date next to the already-existing activity name.
$tasks.Range("B47").Value = 44927
$tasks.Range("B47").NumberFormat = "d-mmm-yy"
$tasks.Range("C47").Value = 44927
$tasks.Range("C47").NumberFormat = "d-mmm-yy"

# Row 48 - new activity with dates (label in column D added further below,
# to line up with the order the strings were originally entered).
$tasks.Range("A48").Value = "end to end "
$tasks.Range("B48").Value = 44940
$tasks.Range("B48").NumberFormat = "d-mmm-yy"
$tasks.Range("C48").Value = 44940
$tasks.Range("C48").NumberFormat = "d-mmm-yy"

# Row 49 intentionally left blank (matches the gap in the original sheet).

# Row 50 - new activity with dates.
$tasks.Range("A50").Value = "Log Writing to File "
$tasks.Range("B50").Value = 44940
$tasks.Range("B50").NumberFormat = "d-mmm-yy"
$tasks.Range("C50").Value = 44940
$tasks.Range("C50").NumberFormat = "d-mmm-yy"

# Rows 51-53 - plain new activity names (no dates).
$tasks.Range("A51").Value = "Converting to a desktop app using enzyme"
$tasks.Range("A52").Value = "Encryption of sensitive data in .env files"
$tasks.Range("A53").Value = "Salary Screen"

# Bold "section" label in column D of row 48.
$tasks.Range("D48").Value = "open bugs"
$tasks.Range("D48").Font.Bold = $true

$tasks.Range("A54").Value = "Doubt : where is leaves availed used "

# ---------------------------------------------------------------------------
# 3. Add the new bug entries to the "Bugs" sheet.
# ---------------------------------------------------------------------------
$bugs = $wb.Worksheets.Item("Bugs")

$bugs.Range("A5").Value = "once Alert Panel is dismissed not able to render new alerts "

$bugs.Range("A6").Value = "No message shown when payroll is not saved..instead shows a stack trac"
$bugs.Range("B6").Value = "fixed"

# ---------------------------------------------------------------------------
# 4. Back to "Tasks" row 55, then "Bugs" row 7 - both use the same new
#    string, entered here for the first time.
# ---------------------------------------------------------------------------
$tasks.Range("A55").Value = "Disable auto complete "
$tasks.Range("B55").Value = 44940
$tasks.Range("B55").NumberFormat = "d-mmm-yy"
$tasks.Range("C55").Value = 44940
$tasks.Range("C55").NumberFormat = "d-mmm-yy"

$bugs.Range("A7").Value = "Disable auto complete "
$bugs.Range("B7").Value = "fixed"

$bugs.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 5. Restore view state: "Tasks" is the active / selected sheet again, with
#    the selection sitting on D52; "Bugs" keeps its selection on A7.
# ---------------------------------------------------------------------------
$bugs.Range("A7").Select() | Out-Null

$tasks.Activate()
$tasks.Range("D52").Select() | Out-Null
